$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column at N (shifts old N:P -> O:Q) on the
# "Repayment schedule" sheet, then size it like its left neighbour (M).
[void]$ws.Columns("N:N").Insert()
$ws.Columns("N:N").ColumnWidth = $ws.Columns("M:M").ColumnWidth

# Make "Repayment schedule" the active sheet/tab, with R6 selected.
[void]$ws.Activate()
[void]$ws.Range("R6").Select()
